# Append a new row (row 81) of logged data to each of the four sheets,
# mirroring the format of the existing rows (row 80) and extending the
# sheet's used range from A1:I80 to A1:I81.

$wb = $excel.ActiveWorkbook

# --- Sheet "FE_LFT_#1" ---
$ws = $wb.Worksheets.Item("FE_LFT_#1")
$ws.Range("A81").Value = 45867.4906712963
$ws.Range("A81").NumberFormat = $ws.Range("A80").NumberFormat
$ws.Range("B81").Value = "0x01,0x7c"
$ws.Range("C81").Value = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
$ws.Range("D81").Value = "0x01,0x20"
$ws.Range("E81").Value = "0xf"
$ws.Range("F81").Value = 380
$ws.Range("G81").Value = 759863127514710945038336.0
$ws.Range("H81").Value = 288
$ws.Range("I81").Value = 15

# --- Sheet "FE_LFT_#2" ---
$ws = $wb.Worksheets.Item("FE_LFT_#2")
$ws.Range("A81").Value = 45867.4906712963
$ws.Range("A81").NumberFormat = $ws.Range("A80").NumberFormat
$ws.Range("B81").Value = "0x01,0x90"
$ws.Range("C81").Value = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
$ws.Range("D81").Value = "0x01,0x2C"
$ws.Range("E81").Value = "0xe"
$ws.Range("F81").Value = 400
$ws.Range("G81").Value = 568432987514711010443264.0
$ws.Range("H81").Value = 300
$ws.Range("I81").Value = 14

# --- Sheet "FE_PLT_#1" ---
$ws = $wb.Worksheets.Item("FE_PLT_#1")
$ws.Range("A81").Value = 45867.4906712963
$ws.Range("A81").NumberFormat = $ws.Range("A80").NumberFormat
$ws.Range("B81").Value = "0x00,0x6e"
$ws.Range("C81").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws.Range("D81").Value = "0x00,0x60"
$ws.Range("E81").Value = "0x3"
$ws.Range("F81").Value = 110
$ws.Range("G81").Value = 568631262647113970876416.0
$ws.Range("H81").Value = 96
$ws.Range("I81").Value = 3

# --- Sheet "FE_PLT_#2" ---
$ws = $wb.Worksheets.Item("FE_PLT_#2")
$ws.Range("A81").Value = 45867.4906712963
$ws.Range("A81").NumberFormat = $ws.Range("A80").NumberFormat
$ws.Range("B81").Value = "0x00,0x6e"
$ws.Range("C81").Value = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
$ws.Range("D81").Value = "0x00,0x5E"
$ws.Range("E81").Value = "0x3"
$ws.Range("F81").Value = 110
$ws.Range("G81").Value = 985046333984776009023488.0
$ws.Range("H81").Value = 94
$ws.Range("I81").Value = 3
